$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 11 of the Rules table (the "R40" rule) now stores the literal text
# "1" in column B instead of the "R40" label. Force the cell to text so
# the new value is written as a string (not auto-converted to a number),
# matching the original cell's text type.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
